# sample-jxls7.xlsx demo report: rename the jxls bean-property placeholders
# used in the "Item" / "Volume" template cells from the old upper-case
# JXLS_VAR style (${row.ITEM_NAME} / ${row.VOLUME}) to the new lower-case
# bean property names (${row.item_name} / ${row.volume}), then leave the
# selection where the author last left it when saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A5 / B5 hold the per-row template placeholders (shared strings) that the
# jxls engine substitutes for each data row under the "Item"/"Volume"
# header in row 4.
$ws.Range("A5").Value = '${row.item_name}'
$ws.Range("B5").Value = '${row.volume}'

# The saved sheet view now has B6 selected instead of A5.
$ws.Range("B6").Select()
